# LOM3089.docx edit: rotate several section bodies' text content.
#
# Net structural effect (verified against the target XML): the number of
# paragraphs, their order, and their paragraph styles are unchanged; only
# the text *content* of specific paragraphs/runs is swapped around. So the
# whole edit can be done with scoped (per-paragraph) Find & Replace, which
# avoids any risk of one replacement text being matched again by a later
# replacement.

$d = $word.ActiveDocument

function Replace-InParagraph([int]$index, [string]$oldText, [string]$newText) {
    $rng = $d.Paragraphs($index).Range
    $ok = $rng.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)
    if (-not $ok) {
        throw "Replace failed in paragraph $index (text not found)"
    }
}

$BR = [char]11

# ---- Objetivos body (paragraph 6): old "Apresentar nocoes..." -> new "Fundamentos..."
$apresentarText = 'Apresentar noções de mecânica dos fluidos e reologia, mediante estudo dos meios fluidos quando estáticos ou em movimento. Capacitar o aluno a modelar e resolver problemas de interesse em mecânica dos fluidos e reologia, com escolha adequada de hipóteses e aplicação de ferramentas correspondentes de solução.'
$fundamentosText = 'Fundamentos de mecânica dos fluidos. Revisão de estática dos fluidos. Formulação integral e diferencial das equações de transporte de massa, energia e quantidade de movimento. Análise dimensional e semelhança. Escoamento incompressível de fluidos ideais e viscosos, regime laminar e turbulento. Equação de Navier-Stokes. Teoria da camada limite. Escoamento de fluidos não newtonianos. Formulação tensorial: tensão e deformação. Viscosidade e reometria. Viscoelasticidade. Aplicações.'
Replace-InParagraph 6 $apresentarText $fundamentosText

# ---- Docente(s) list item (paragraph 8): old "519033 - Carlos..." -> new "Apresentar nocoes..."
$carlosText = '519033 - Carlos Yujiro Shigue'
Replace-InParagraph 8 $carlosText $apresentarText

# ---- Programa resumido body (paragraph 10): old "Fundamentos..." -> new "Introducao: conceito de fluido..."
$introducaoText = 'Introdução: conceito de fluido; propriedades e conceito de contínuo; modelagem de processos de transferência; métodos de análise; dimensões e unidades.' + $BR + `
    'Revisão de estática de fluidos: equação básica da hidrostática, variação de pressão em um fluido estático; princípios de Stevin, de Pascal e de Arquimedes.' + $BR + `
    'Formulação integral das equações de transporte: teorema de transporte de Reynolds; aplicação para os princípios de conservação de massa, quantidade de movimento e energia; equação de Bernoulli.' + $BR + `
    'Formulação diferencial das equações de transporte: descrição do escoamento; forma diferencial: dos princípios de conservação de massa, quantidade de movimento e energia; formulação adimensional, análise dimensional e semelhança. Grupos adimensionais: número de Reynolds e número de Grashoff.' + $BR + `
    'Escoamento incompressível interno: equações de Euler; lei de Newton para a viscosidade, tensões de cisalhamento; equação de Navier-Stokes; regimes de escoamento: escoamento laminar e turbulento. Cálculo de perda de carga (distribuída e localizada), coeficiente de atrito. ' + $BR + `
    'Escoamento incompressível externo: introdução à camada limite; escoamento ao redor de corpos, força da arraste.' + $BR + `
    'Introdução a reologia. Definição e formulação tensorial de tensão e deformação. Tipos de deformação e escoamento de materiais. Equações fundamentais da reologia. Escoamento de fluidos newtonianos e não newtonianos. Viscosimetria e reometria. Reologia de sistemas dispersos. Colóides e emulsões. Soluções diluídas. Viscosimetria capilar. Aplicações.'
Replace-InParagraph 10 $fundamentosText $introducaoText

# ---- Avaliacao body (paragraph 12): old "Introducao: conceito de fluido..." -> new "A avaliacao sera feita..."
$avaliacaoText = 'A avaliação será feita por meio de duas provas escritas P1 e P2 e por listas de exercícios e relatórios.'
Replace-InParagraph 12 $introducaoText $avaliacaoText

# ---- Avaliacao criteria paragraph (paragraph 14): three labelled values shift down one slot.
# NOTE: cannot chain simple scoped Find/Replace calls here, because the new
# text of one value can equal the pre-existing text of a sibling value
# (e.g. "Metodo:" becomes the old "Critério:" text), which would make a
# later ReplaceAll match more than one run. Instead, locate each label's
# own range precisely (labels are unique) and set the text of the value
# that immediately follows it, bounded by the following label (or the end
# of the paragraph).
$notaFinalText = 'A Nota final (NF) será calculada pela média ponderada das provas escritas e pela média dos trabalhos TR da seguinte maneira: NF = (P1 + 2*P2 + TR)/4'
$recuperacaoText = 'A recuperação será feita por meio de uma prova escrita (PR) e a média de recuperação (MR) calculada pela fórmula: MR = (NF + PR)/2'
$bibliografiaText = 'BIRD,R. B.; STEWART, W. E.; LIGHTFOOT, E. N. Fenômenos de Transporte. LTC Editora, 2004.' + $BR + `
    'FOX, R. W., McDONALD, A. T. Introdução à Mecânica dos Fluidos. LTC Editora, 2001.' + $BR + `
    'SISSOM, L. E., PITTS, D. R. Fenômenos de Transporte. Ed. Guanabara, 1988.' + $BR + `
    'SCHRAMM, G. Reologia e Reometria. Editora Artliber, 2006.' + $BR + `
    'MANRICH, S.; PESSAN, L.A. Reologia: Conceitos Básicos, Editora UFSCar, 1987.' + $BR + `
    'MALKIN, A. Rheology Fundamentals. ChemTec Publishing, 1994.'

function Find-LabelRange([int]$searchStart, [int]$searchEnd, [string]$label) {
    $r = $d.Range($searchStart, $searchEnd)
    $ok = $r.Find.Execute($label, $true, $false, $false, $false, $false, $true, 0, $false, "", 0)
    if (-not $ok) {
        throw "label not found: $label"
    }
    return $r
}

$p14Range = $d.Paragraphs(14).Range
$p14Start = $p14Range.Start
$p14End = $p14Range.End - 1   # exclude the paragraph mark

$metodoLabel = Find-LabelRange $p14Start $p14End 'Método: '
$criterioLabel = Find-LabelRange $metodoLabel.End $p14End 'Critério: '
$normaLabel = Find-LabelRange $criterioLabel.End $p14End 'Norma de recuperação: '

# Replace in reverse order so earlier ranges' offsets stay valid while later
# (higher-offset) ranges are still being rewritten.
$d.Range($normaLabel.End, $p14End).Text = $bibliografiaText
$d.Range($criterioLabel.End, $normaLabel.Start).Text = $recuperacaoText + $BR
$d.Range($metodoLabel.End, $criterioLabel.Start).Text = $notaFinalText + $BR

# ---- Bibliografia body (paragraph 16): old bibliography list -> new "519033 - Carlos..."
Replace-InParagraph 16 $bibliografiaText $carlosText
